# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" everywhere it
#   appears (Overview summary columns + per-language Status column).
# - Handoff/handback timestamps bump forward a few minutes to reflect the
#   new report generation time.
# - The "handback file is not the latest" message now points at the newest
#   commit hash.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$overview.Range("G2").Value = "2017-02-10 07:33:41"
$dede.Range("H2").Value = "2017-02-10 07:33:41"
$zhcn.Range("H2").Value = "2017-02-10 07:33:23"

# --- Latest handback-file-is-stale message: bump the "latest" commit hash ---
$newUrl = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f665d099d04b5d8a45f082eb5335899f77b2133/e2e/68e93543-7085-44ff-8670-1b4d66bd4f4f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7ffe4e366b6f35e1c97379bcbe7cc496a49ea500/e2e/68e93543-7085-44ff-8670-1b4d66bd4f4f.md."

$zhcn.Range("R2").Value = $newUrl
$dede.Range("R2").Value = $newUrl

# --- Column widths: widen the "Status" columns now that "Ready for handoff"
#     is longer than "In Translation" (target stored width ~17.216 chars;
#     16.3 is the closest ColumnWidth input that round-trips to it) ---
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
